$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.895.31"
$ws.Range("E2").Value = "  +3.34%  "
$ws.Range("D3").Value = "3.135.52"
$ws.Range("E3").Value = "  +2.48%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.63"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.15"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.127.83"
$ws.Range("E8").Value = "  +2.58%  "
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +19.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.73"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.12"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.87%  "
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "3.650.73"
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "63.818.97"
$ws.Range("E18").Value = "  +3.43%  "
$ws.Range("D19").Value = "3.132.34"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.62"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.26"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.58"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.34"
$ws.Range("D24").ClearFormats()
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.68"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +8.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.72"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.23%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.87"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.09"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("E34").Value = "  +7.96%  "
$ws.Range("E35").Value = "  +9.45%  "
$ws.Range("E36").Value = "  +2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.42"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +16.09%  "
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "456.15"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +9.81%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.88"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.73"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("D43").Value = "2.921.75"
$ws.Range("E43").Value = "  +5.42%  "
$ws.Range("E44").Value = "  +5.18%  "
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "129.60"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.33%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.50"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.61%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.74"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.85%  "
